# Update loading_percent values for Case_2_148 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 12.99114768829283
$ws.Range("C2").Value = 6.443580603189078
$ws.Range("E2").Value = 11.5266319354333
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 41.17971849164739
$ws.Range("H2").Value = 17.49628690231036
$ws.Range("I2").Value = 27.41399710471474
$ws.Range("K2").Value = 10.5933071817574
$ws.Range("L2").Value = 10.24090758299286
$ws.Range("M2").Value = 14.61703621086329
$ws.Range("N2").Value = 21.14889742969934

# Row 3
$ws.Range("B3").Value = 12.81077756571739
$ws.Range("C3").Value = 6.362136940330569
$ws.Range("E3").Value = 11.54635266528762
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 41.28840621595985
$ws.Range("H3").Value = 17.54965785045948
$ws.Range("I3").Value = 27.50768331018679
$ws.Range("K3").Value = 10.46868518184184
$ws.Range("L3").Value = 10.24947439209668
$ws.Range("M3").Value = 14.59650511187795
$ws.Range("N3").Value = 21.21399276139048

# Row 4
$ws.Range("B4").Value = 12.70146751636801
$ws.Range("C4").Value = 6.310731392090664
$ws.Range("E4").Value = 11.55976487623375
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 41.36570089473734
$ws.Range("H4").Value = 17.58504048248533
$ws.Range("I4").Value = 27.5697913099075
$ws.Range("K4").Value = 10.39354870916268
$ws.Range("L4").Value = 10.25613126960129
$ws.Range("M4").Value = 14.58612044560316
$ws.Range("N4").Value = 21.25582921700659

# Row 5
$ws.Range("B5").Value = 12.65733921835183
$ws.Range("C5").Value = 6.289440569519396
$ws.Range("E5").Value = 11.56555867852502
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 41.39984580248006
$ws.Range("H5").Value = 17.60011612855523
$ws.Range("I5").Value = 27.59625296457683
$ws.Range("K5").Value = 10.36331195090448
$ws.Range("L5").Value = 10.25919574146181
$ws.Range("M5").Value = 14.58245091455587
$ws.Range("N5").Value = 21.27334893769841

# Row 6
$ws.Range("B6").Value = 12.65003849217224
$ws.Range("C6").Value = 6.285884820465183
$ws.Range("E6").Value = 11.5665405704292
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 41.40567512677261
$ws.Range("H6").Value = 17.60265910444907
$ws.Range("I6").Value = 27.60071647078741
$ws.Range("K6").Value = 10.35831522376082
$ws.Range("L6").Value = 10.25972585218963
$ws.Range("M6").Value = 14.58187565475591
$ws.Range("N6").Value = 21.27628656355494

# Row 7
$ws.Range("B7").Value = 12.70087062850476
$ws.Range("C7").Value = 6.310445631758061
$ws.Range("E7").Value = 11.55984168399003
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 41.36615068048555
$ws.Range("H7").Value = 17.58524113806697
$ws.Range("I7").Value = 27.57014351688683
$ws.Range("K7").Value = 10.39313933408087
$ws.Range("L7").Value = 10.25617117334628
$ws.Range("M7").Value = 14.58606867565207
$ws.Range("N7").Value = 21.25606358509007

# Row 8
$ws.Range("B8").Value = 12.92869066870381
$ws.Range("C8").Value = 6.41579452461543
$ws.Range("E8").Value = 11.53316133369646
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 41.21499724311777
$ws.Range("H8").Value = 17.51414692782566
$ws.Range("I8").Value = 27.44534826974522
$ws.Range("K8").Value = 10.55007156145073
$ws.Range("L8").Value = 10.24357182488981
$ws.Range("M8").Value = 14.60949832098763
$ws.Range("N8").Value = 21.17095550667028

# Row 9
$ws.Range("B9").Value = 13.38443067904657
$ws.Range("C9").Value = 6.61090067609443
$ws.Range("E9").Value = 11.49116624568728
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 41.00273895826813
$ws.Range("H9").Value = 17.39546328890172
$ws.Range("I9").Value = 27.23702697254018
$ws.Range("K9").Value = 10.86724008947148
$ws.Range("L9").Value = 10.22992408943339
$ws.Range("M9").Value = 14.67289575655889
$ws.Range("N9").Value = 21.01881341249116

# Row 10
$ws.Range("B10").Value = 13.72150613525355
$ws.Range("C10").Value = 6.746754702935307
$ws.Range("E10").Value = 11.46658158607352
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 40.89855623131024
$ws.Range("H10").Value = 17.32090737934566
$ws.Range("I10").Value = 27.10620237731714
$ws.Range("K10").Value = 11.10392790315084
$ws.Range("L10").Value = 10.22660496804245
$ws.Range("M10").Value = 14.72986041970407
$ws.Range("N10").Value = 20.91593758114567

# Row 11
$ws.Range("B11").Value = 13.87466501518176
$ws.Range("C11").Value = 6.80683402392047
$ws.Range("E11").Value = 11.45675302296595
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 40.86248523978638
$ws.Range("H11").Value = 17.28973545058166
$ws.Range("I11").Value = 27.05152202629767
$ws.Range("K11").Value = 11.21195784752339
$ws.Range("L11").Value = 10.22654284714556
$ws.Range("M11").Value = 14.75796992026403
$ws.Range("N11").Value = 20.87105028190102

# Row 12
$ws.Range("B12").Value = 13.9325821112037
$ws.Range("C12").Value = 6.829329114633185
$ws.Range("E12").Value = 11.45322557335095
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 40.85045986612125
$ws.Range("H12").Value = 17.27832613285059
$ws.Range("I12").Value = 27.0315116815575
$ws.Range("K12").Value = 11.2528812942433
$ws.Range("L12").Value = 10.22672662680449
$ws.Range("M12").Value = 14.76892435685579
$ws.Range("N12").Value = 20.85432604346981

# Row 13
$ws.Range("B13").Value = 13.92011308496632
$ws.Range("C13").Value = 6.824495889715195
$ws.Range("E13").Value = 11.45397663351029
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 40.85297699248282
$ws.Range("H13").Value = 17.28076577019801
$ws.Range("I13").Value = 27.0357902999083
$ws.Range("K13").Value = 11.24406761631709
$ws.Range("L13").Value = 10.22667784159154
$ws.Range("M13").Value = 14.76655143224431
$ws.Range("N13").Value = 20.85791575852768

# Row 14
$ws.Range("B14").Value = 13.87943182003246
$ws.Range("C14").Value = 6.808689877307185
$ws.Range("E14").Value = 11.45645892386731
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 40.8614631283968
$ws.Range("H14").Value = 17.28878888507588
$ws.Range("I14").Value = 27.04986181010689
$ws.Range("K14").Value = 11.21532455499365
$ws.Range("L14").Value = 10.2265538174769
$ws.Range("M14").Value = 14.75886497026051
$ws.Range("N14").Value = 20.86966889237972

# Row 15
$ws.Range("B15").Value = 13.85450123215733
$ws.Range("C15").Value = 6.798974710982732
$ws.Range("E15").Value = 11.4580047031168
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 40.86687407551665
$ws.Range("H15").Value = 17.2937546983699
$ws.Range("I15").Value = 27.05857167578688
$ws.Range("K15").Value = 11.19771945255054
$ws.Range("L15").Value = 10.22650481860262
$ws.Range("M15").Value = 14.75419698121292
$ws.Range("N15").Value = 20.87690362098712

# Row 16
$ws.Range("B16").Value = 13.71148849701764
$ws.Range("C16").Value = 6.742793038290943
$ws.Range("E16").Value = 11.46725113618075
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 40.90114190270447
$ws.Range("H16").Value = 17.32299977369257
$ws.Range("I16").Value = 27.1098731946952
$ws.Range("K16").Value = 11.09687196127123
$ws.Range("L16").Value = 10.2266380868029
$ws.Range("M16").Value = 14.72806709788919
$ws.Range("N16").Value = 20.91890942716848

# Row 17
$ws.Range("B17").Value = 13.62366934382989
$ws.Range("C17").Value = 6.707880637392351
$ws.Range("E17").Value = 11.47327028761736
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 40.92506839147008
$ws.Range("H17").Value = 17.34164363773855
$ws.Range("I17").Value = 27.14258343842961
$ws.Range("K17").Value = 11.03507094455544
$ws.Range("L17").Value = 10.22709004621484
$ws.Range("M17").Value = 14.71259567881334
$ws.Range("N17").Value = 20.94516729658696

# Row 18
$ws.Range("B18").Value = 13.57314450280922
$ws.Range("C18").Value = 6.687638389197205
$ws.Range("E18").Value = 11.47685992295535
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 40.93989568276296
$ws.Range("H18").Value = 17.35262535219981
$ws.Range("I18").Value = 27.16185229369691
$ws.Range("K18").Value = 10.99956063783473
$ws.Range("L18").Value = 10.22748632953773
$ws.Range("M18").Value = 14.70390399217931
$ws.Range("N18").Value = 20.96045011586958

# Row 19
$ws.Range("B19").Value = 13.55603703407694
$ws.Range("C19").Value = 6.680757216844501
$ws.Range("E19").Value = 11.47809723698708
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 40.94509874358066
$ws.Range("H19").Value = 17.35638792057621
$ws.Range("I19").Value = 27.16845447624332
$ws.Range("K19").Value = 10.98754477502314
$ws.Range("L19").Value = 10.22764394438663
$ws.Range("M19").Value = 14.70099687155673
$ws.Range("N19").Value = 20.96565556728342

# Row 20
$ws.Range("B20").Value = 13.63301965373471
$ws.Range("C20").Value = 6.711613902668604
$ws.Range("E20").Value = 11.4726163379748
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 40.92241106553602
$ws.Range("H20").Value = 17.33963223619497
$ws.Range("I20").Value = 27.1390543008542
$ws.Range("K20").Value = 11.04164631105481
$ws.Range("L20").Value = 10.22702783105406
$ws.Range("M20").Value = 14.71422124842361
$ws.Range("N20").Value = 20.94235348265959

# Row 21
$ws.Range("B21").Value = 13.89138351942501
$ws.Range("C21").Value = 6.813339487353161
$ws.Range("E21").Value = 11.45572454236677
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 40.85892615806949
$ws.Range("H21").Value = 17.2864215846541
$ws.Range("I21").Value = 27.04570977130323
$ws.Range("K21").Value = 11.22376697590313
$ws.Range("L21").Value = 10.22658462763851
$ws.Range("M21").Value = 14.76111430551792
$ws.Range("N21").Value = 20.86620929637347

# Row 22
$ws.Range("B22").Value = 14.05974362603072
$ws.Range("C22").Value = 6.878329430109837
$ws.Range("E22").Value = 11.44581774807313
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 40.8269609894449
$ws.Range("H22").Value = 17.25394674517634
$ws.Range("I22").Value = 26.98876054044409
$ws.Range("K22").Value = 11.34286336785624
$ws.Range("L22").Value = 10.22750282652781
$ws.Range("M22").Value = 14.79356551368058
$ws.Range("N22").Value = 20.81803888941221

# Row 23
$ws.Range("B23").Value = 13.96994980806151
$ws.Range("C23").Value = 6.843782330353696
$ws.Range("E23").Value = 11.45100167757495
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 40.84314807899791
$ws.Range("H23").Value = 17.27106852026479
$ws.Range("I23").Value = 27.01878387961039
$ws.Range("K23").Value = 11.27930495727637
$ws.Range("L23").Value = 10.22690255959044
$ws.Range("M23").Value = 14.77608263439035
$ws.Range("N23").Value = 20.84360288235885

# Row 24
$ws.Range("B24").Value = 13.62879248789472
$ws.Range("C24").Value = 6.709926624370957
$ws.Range("E24").Value = 11.47291158631014
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 40.92360910493658
$ws.Range("H24").Value = 17.34054077126584
$ws.Range("I24").Value = 27.14064838085442
$ws.Range("K24").Value = 11.03867352298105
$ws.Range("L24").Value = 10.22705553345001
$ws.Range("M24").Value = 14.71348569551915
$ws.Range("N24").Value = 20.94362502588685

# Row 25
$ws.Range("B25").Value = 13.26051504993087
$ws.Range("C25").Value = 6.559400373215885
$ws.Range("E25").Value = 11.50142421334424
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 41.05110087602046
$ws.Range("H25").Value = 17.42535064768152
$ws.Range("I25").Value = 27.28948249227768
$ws.Range("K25").Value = 10.7806433512379
$ws.Range("L25").Value = 10.23243541422358
$ws.Range("M25").Value = 14.65390145953185
$ws.Range("N25").Value = 21.0584018329046
